# Apply updates described by the diff:
# Sheet "All Orders" (sheet1): row 19 - Status -> CANCELLED, Cancel Reason -> "test order"
# Sheet "Daily Summary" (sheet2): row 4 - Cancelled 5->6, Revenue 155->125, Pending 155->125

$wb = $excel.ActiveWorkbook

$wsOrders = $wb.Worksheets.Item("All Orders")
$wsOrders.Range("H19").Value = "CANCELLED"
$wsOrders.Range("M19").Value = "test order"

$wsSummary = $wb.Worksheets.Item("Daily Summary")
$wsSummary.Range("D4").Value = 6
$wsSummary.Range("E4").Value = 125
$wsSummary.Range("G4").Value = 125
